$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.86352615203185
$ws.Range("C2").Value = 11.0704057494741
$ws.Range("E2").Value = 15.80939750668612
$ws.Range("F2").Value = 37.05606187625315
$ws.Range("G2").Value = 3.648498030116251
$ws.Range("I2").Value = 24.07044985236733
$ws.Range("J2").Value = 7.8939449450746
$ws.Range("L2").Value = 12.73870787281835
$ws.Range("M2").Value = 16.6924763437016
$ws.Range("O2").Value = 23.6314909350026
$ws.Range("B3").Value = 15.31927779028043
$ws.Range("C3").Value = 10.87355907233613
$ws.Range("E3").Value = 15.84964383976237
$ws.Range("F3").Value = 37.16121959947497
$ws.Range("G3").Value = 3.650514273403191
$ws.Range("I3").Value = 24.21619582048301
$ws.Range("J3").Value = 7.886602316594892
$ws.Range("L3").Value = 12.72932254275976
$ws.Range("M3").Value = 16.55784133550932
$ws.Range("O3").Value = 23.74433372204495
$ws.Range("B4").Value = 14.97588008924854
$ws.Range("C4").Value = 10.75015200977517
$ws.Range("E4").Value = 15.87599280447426
$ws.Range("F4").Value = 37.2351785390731
$ws.Range("G4").Value = 3.651817802661525
$ws.Range("I4").Value = 24.31100965874567
$ws.Range("J4").Value = 7.882265564305893
$ws.Range("L4").Value = 12.72494844834581
$ws.Range("M4").Value = 16.47622150682558
$ws.Range("O4").Value = 23.81969192153575
$ws.Range("B5").Value = 14.83381850979123
$ws.Range("C5").Value = 10.69926627003389
$ws.Range("E5").Value = 15.88714276165102
$ws.Range("F5").Value = 37.26767269636635
$ws.Range("G5").Value = 3.652365535868287
$ws.Range("I5").Value = 24.35098663475663
$ws.Range("J5").Value = 7.880542012092248
$ws.Range("L5").Value = 12.72351700587047
$ws.Range("M5").Value = 16.44324961063314
$ws.Range("O5").Value = 23.85192371695996
$ws.Range("B6").Value = 14.81010708643467
$ws.Range("C6").Value = 10.69078199923114
$ws.Range("E6").Value = 15.88901914247477
$ws.Range("F6").Value = 37.2732103859766
$ws.Range("G6").Value = 3.652457486717867
$ws.Range("I6").Value = 24.35770570373376
$ws.Range("J6").Value = 7.880258473336397
$ws.Range("L6").Value = 12.72330057011306
$ws.Range("M6").Value = 16.43779283826296
$ws.Range("O6").Value = 23.85736762720773
$ws.Range("B7").Value = 14.97397252543015
$ws.Range("C7").Value = 10.74946810361651
$ws.Range("E7").Value = 15.87614150531367
$ws.Range("F7").Value = 37.23560723889177
$ws.Range("G7").Value = 3.651825122566013
$ws.Range("I7").Value = 24.31154337767583
$ws.Range("J7").Value = 7.88224214218823
$ws.Range("L7").Value = 12.7249277196876
$ws.Range("M7").Value = 16.47577563324333
$ws.Range("O7").Value = 23.82012045096779
$ws.Range("B8").Value = 15.6778962321117
$ws.Range("C8").Value = 11.00308594919989
$ws.Range("E8").Value = 15.82293502802693
$ws.Range("F8").Value = 37.09036735558091
$ws.Range("G8").Value = 3.649179657695894
$ws.Range("I8").Value = 24.11959810057782
$ws.Range("J8").Value = 7.891377908594515
$ws.Range("L8").Value = 12.7351847802825
$ws.Range("M8").Value = 16.64585250614599
$ws.Range("O8").Value = 23.66913630810757
$ws.Range("B9").Value = 16.97739902000403
$ws.Range("C9").Value = 11.47845790008043
$ws.Range("E9").Value = 15.73155606884523
$ws.Range("F9").Value = 36.88033465257362
$ws.Range("G9").Value = 3.644509628061271
$ws.Range("I9").Value = 23.78542468817365
$ws.Range("J9").Value = 7.910631663868174
$ws.Range("L9").Value = 12.76623019004021
$ws.Range("M9").Value = 16.98653162141055
$ws.Range("O9").Value = 23.42142454510838
$ws.Range("B10").Value = 17.87380370998409
$ws.Range("C10").Value = 11.81204229853064
$ws.Range("E10").Value = 15.67227228029373
$ws.Range("F10").Value = 36.77194407796073
$ws.Range("G10").Value = 3.641390847282964
$ws.Range("I10").Value = 23.56562127105209
$ws.Range("J10").Value = 7.92556279341809
$ws.Range("L10").Value = 12.79558259885839
$ws.Range("M10").Value = 17.23964022463923
$ws.Range("O10").Value = 23.26917105063877
$ws.Range("B11").Value = 18.26736076480914
$ws.Range("C11").Value = 11.95995279750836
$ws.Range("E11").Value = 15.64699755374953
$ws.Range("F11").Value = 36.73266514311375
$ws.Range("G11").Value = 3.640039136295371
$ws.Range("I11").Value = 23.47120562921459
$ws.Range("J11").Value = 7.932519268343939
$ws.Range("L11").Value = 12.8103283597408
$ws.Range("M11").Value = 17.35506293530111
$ws.Range("O11").Value = 23.20641889975172
$ws.Range("B12").Value = 18.41423478759492
$ws.Range("C12").Value = 12.01537633783902
$ws.Range("E12").Value = 15.63766947747003
$ws.Range("F12").Value = 36.71923718881533
$ws.Range("G12").Value = 3.639536864921527
$ws.Range("I12").Value = 23.43625442506615
$ws.Range("J12").Value = 7.935176518384576
$ws.Range("L12").Value = 12.81610989309999
$ws.Range("M12").Value = 17.39878301354151
$ws.Range("O12").Value = 23.18359672574833
$ws.Range("B13").Value = 18.38270043669711
$ws.Range("C13").Value = 12.00346652560479
$ws.Range("E13").Value = 15.63966765216295
$ws.Range("F13").Value = 36.72206476478219
$ws.Range("G13").Value = 3.639644612187629
$ws.Range("I13").Value = 23.44374611762705
$ws.Range("J13").Value = 7.934603219983151
$ws.Range("L13").Value = 12.81485599110241
$ws.Range("M13").Value = 17.38936703204311
$ws.Range("O13").Value = 23.18846997596249
$ws.Range("B14").Value = 18.27948801297553
$ws.Range("C14").Value = 11.9645244816044
$ws.Range("E14").Value = 15.64622526262156
$ws.Range("F14").Value = 36.73153141388204
$ws.Range("G14").Value = 3.639997622156267
$ws.Range("I14").Value = 23.46831409761301
$ws.Range("J14").Value = 7.932737424923346
$ws.Range("L14").Value = 12.81080006505684
$ws.Range("M14").Value = 17.35865973507336
$ws.Range("O14").Value = 23.2045224255783
$ws.Range("B15").Value = 18.21598335551474
$ws.Range("C15").Value = 11.94059391538048
$ws.Range("E15").Value = 15.65027360387203
$ws.Range("F15").Value = 36.73751844691993
$ws.Range("G15").Value = 3.640215098737668
$ws.Range("I15").Value = 23.48346714630897
$ws.Range("J15").Value = 7.931597545855635
$ws.Range("L15").Value = 12.80834134845397
$ws.Range("M15").Value = 17.3398513463321
$ws.Range("O15").Value = 23.21447766224241
$ws.Range("B16").Value = 17.84778772061698
$ws.Range("C16").Value = 11.80229571822732
$ws.Range("E16").Value = 15.67395807051997
$ws.Range("F16").Value = 36.77471320778861
$ws.Range("G16").Value = 3.641480529661523
$ws.Range("I16").Value = 23.57190375188906
$ws.Range("J16").Value = 7.92511142231266
$ws.Range("L16").Value = 12.79464673503874
$ws.Range("M16").Value = 17.23210031229782
$ws.Range("O16").Value = 23.27340343113256
$ws.Range("B17").Value = 17.61818457755011
$ws.Range("C17").Value = 11.71644552113417
$ws.Range("E17").Value = 15.68892107495967
$ws.Range("F17").Value = 36.80010260489171
$ws.Range("G17").Value = 3.642273966504978
$ws.Range("I17").Value = 23.62758464875501
$ws.Range("J17").Value = 7.921174011838573
$ws.Range("L17").Value = 12.7866004808425
$ws.Range("M17").Value = 17.16605046440856
$ws.Range("O17").Value = 23.31122292117611
$ws.Range("B18").Value = 17.48479007035962
$ws.Range("C18").Value = 11.66670739554633
$ws.Range("E18").Value = 15.69768686292266
$ws.Range("F18").Value = 36.81564964128016
$ws.Range("G18").Value = 3.642736643158459
$ws.Range("I18").Value = 23.66013545142217
$ws.Range("J18").Value = 7.918924813859786
$ws.Range("L18").Value = 12.78210379330886
$ws.Range("M18").Value = 17.12808861854248
$ws.Range("O18").Value = 23.33358804182985
$ws.Range("B19").Value = 17.43939971948488
$ws.Range("C19").Value = 11.64980629983354
$ws.Range("E19").Value = 15.70068221596982
$ws.Range("F19").Value = 36.82107555289213
$ws.Range("G19").Value = 3.642894383208431
$ws.Range("I19").Value = 23.67124668770049
$ws.Range("J19").Value = 7.918165954704591
$ws.Range("L19").Value = 12.78060392431553
$ws.Range("M19").Value = 17.11524112065206
$ws.Range("O19").Value = 23.34126551353195
$ws.Range("B20").Value = 17.64276502779464
$ws.Range("C20").Value = 11.72562187383949
$ws.Range("E20").Value = 15.68731173842813
$ws.Range("F20").Value = 36.79730215669222
$ws.Range("G20").Value = 3.642188850797324
$ws.Range("I20").Value = 23.62160302299279
$ws.Range("J20").Value = 7.921591556632815
$ws.Range("L20").Value = 12.78744344796025
$ws.Range("M20").Value = 17.17307886428042
$ws.Range("O20").Value = 23.30713356132117
$ws.Range("B21").Value = 18.30986338451255
$ws.Range("C21").Value = 11.97597890166911
$ws.Range("E21").Value = 15.6442925459767
$ws.Range("F21").Value = 36.72871155314687
$ws.Range("G21").Value = 3.639893674619236
$ws.Range("I21").Value = 23.46107611786225
$ws.Range("J21").Value = 7.933284835462419
$ws.Range("L21").Value = 12.81198604734113
$ws.Range("M21").Value = 17.36767911382992
$ws.Range("O21").Value = 23.1997818639016
$ws.Range("B22").Value = 18.73323260200879
$ws.Range("C22").Value = 12.13616575979494
$ws.Range("E22").Value = 15.61759261085568
$ws.Range("F22").Value = 36.69231416659028
$ws.Range("G22").Value = 3.638449533336811
$ws.Range("I22").Value = 23.36083715080648
$ws.Range("J22").Value = 7.941060737797912
$ws.Range("L22").Value = 12.82917666974289
$ws.Range("M22").Value = 17.49491846870701
$ws.Range("O22").Value = 23.13510663744752
$ws.Range("B23").Value = 18.50846064738608
$ws.Range("C23").Value = 12.05099640226658
$ws.Range("E23").Value = 15.63171355384985
$ws.Range("F23").Value = 36.71096757627559
$ws.Range("G23").Value = 3.639215200547753
$ws.Range("I23").Value = 23.41390868088627
$ws.Range("J23").Value = 7.936898572788737
$ws.Range("L23").Value = 12.81989735703742
$ws.Range("M23").Value = 17.42701262634227
$ws.Range("O23").Value = 23.16912155718726
$ws.Range("B24").Value = 17.63165653520999
$ws.Range("C24").Value = 11.72147443106513
$ws.Range("E24").Value = 15.68803881058837
$ws.Range("F24").Value = 36.79856527940245
$ws.Range("G24").Value = 3.642227311298089
$ws.Range("I24").Value = 23.62430563668915
$ws.Range("J24").Value = 7.921402739497255
$ws.Range("L24").Value = 12.78706193984763
$ws.Range("M24").Value = 17.16990128796143
$ws.Range("O24").Value = 23.30898042315115
$ws.Range("B25").Value = 16.63550760343065
$ws.Range("C25").Value = 11.35245202360345
$ws.Range("E25").Value = 15.7548941279641
$ws.Range("F25").Value = 36.92911171315973
$ws.Range("G25").Value = 3.645717914546788
$ws.Range("I25").Value = 23.78542468817365
$ws.Range("J25").Value = 7.905282837467609
$ws.Range("L25").Value = 12.76623019004021
$ws.Range("M25").Value = 16.98653162141055
$ws.Range("O25").Value = 23.42142454510838
